$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster data (Player, Position, Team) replacing the previous 18 rows
# with 17 rows (Rui Hachimura and Khris Middleton removed, Nicolas Claxton
# added, remaining players re-sorted).
$data = @(
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Nicolas Claxton", "C", "Brooklyn Nets"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Zion Williamson", "PF,C", "New Orleans Pelicans"),
    @("Kyle Kuzma", "PF", "Washington Wizards")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# The old sheet had 18 data rows (through row 19); the new one has 17
# (through row 18), so remove the now-unused trailing row.
$ws.Rows.Item(19).Delete()
